$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend column A (row counter, styled like existing rows) down to row 34
$ws.Range("A25").Copy()
$ws.Range("A26:A34").PasteSpecial(-4122)
for ($r = 26; $r -le 34; $r++) {
  $ws.Cells.Item($r, 1).Value = $r - 2
}

# Update stock-ticker cells per column for rows 2-34
$ws.Range("B2").Value = "NSE:BALRAMCHIN"
$ws.Range("D2").Value = "NSE:NESTLEIND"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "NSE:GLENMARK"
$ws.Range("B3").Value = "NSE:DPWIRES"
$ws.Range("C3").Value = "NSE:ABBOTINDIA"
$ws.Range("F3").Value = "NSE:HINDPETRO"
$ws.Range("B4").Value = "NSE:HDFCBSE500"
$ws.Range("C4").Value = "NSE:APEX"
$ws.Range("B5").Value = "NSE:HINDPETRO"
$ws.Range("C5").Value = "NSE:ARIHANTCAP"
$ws.Range("B6").Value = "NSE:HINDUNILVR"
$ws.Range("C6").Value = "NSE:ARVEE"
$ws.Range("B7").Value = "NSE:IRIS"
$ws.Range("C7").Value = "NSE:ASHOKLEY"
$ws.Range("B8").Value = "NSE:MASPTOP50"
$ws.Range("C8").Value = "NSE:BAFNAPH"
$ws.Range("B9").Value = "NSE:MAZDA"
$ws.Range("C9").Value = "NSE:BANCOINDIA"
$ws.Range("B10").Value = "NSE:METROPOLIS"
$ws.Range("C10").Value = "NSE:BANKA"
$ws.Range("B11").Value = "NSE:MOHEALTH"
$ws.Range("C11").Value = "NSE:BLUEDART"
$ws.Range("B12").Value = "NSE:MURUDCERA"
$ws.Range("C12").Value = "NSE:CENTRALBK"
$ws.Range("B13").Value = "NSE:PCJEWELLER"
$ws.Range("C13").Value = "NSE:CHEVIOT"
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = "NSE:DANGEE"
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = "NSE:DMART"
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = "NSE:DSSL"
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = "NSE:E2E"
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = "NSE:GHCLTEXTIL"
$ws.Range("B19").Value = ""
$ws.Range("C19").Value = "NSE:GUJGASLTD"
$ws.Range("B20").Value = ""
$ws.Range("C20").Value = "NSE:HINDCOPPER"
$ws.Range("B21").Value = ""
$ws.Range("C21").Value = "NSE:IFCI"
$ws.Range("C22").Value = "NSE:IPL"
$ws.Range("C23").Value = "NSE:JINDALSAW"
$ws.Range("C24").Value = "NSE:LMW"
$ws.Range("C25").Value = "NSE:MANAKSIA"
$ws.Range("C26").Value = "NSE:MRF"
$ws.Range("C27").Value = "NSE:MTNL"
$ws.Range("C28").Value = "NSE:OBEROIRLTY"
$ws.Range("C29").Value = "NSE:PASUPTAC"
$ws.Range("C30").Value = "NSE:RAJRATAN"
$ws.Range("C31").Value = "NSE:RAYMOND"
$ws.Range("C32").Value = "NSE:RELIGARE"
$ws.Range("C33").Value = "NSE:RITES"
$ws.Range("C34").Value = "NSE:SAIL"

Write-Host "done"